# "Data through 6th Ipam" - fill in Day 15 (row 19) and Day 16 (row 20)
# entries for the stress-test schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D (Projected Temp setting label) - same label for both rows
$ws.Range("D19").Value = "30-32"
$ws.Range("D20").Value = "30-32"

# Column G (Current/Projected Setting)
$ws.Range("G19").Value = "87 @ 15:30"
$ws.Range("G20").Value = 87

# Column L (Change Filter)
$ws.Range("L19").Value = "Switch 5 micron to 20 micron"
$ws.Range("L20").Value = "Remove 1 micron filter @ 12:30"

# Column H (Heaters)
$ws.Range("H19").Value = "N"
$ws.Range("H20").Value = "Y @ 12:30 set to 87"

# Move the on-screen selection, like the author's saved view.
[void]$ws.Range("F22").Select()
